{"js": "// Update the worksheet date and every three-digit \u00d7 one-digit answer cell\n// to the new values, per the commit's regenerated output.\nconst replacements = [\n  [\"2024-07-25 Thursday\", \"2024-07-26 Friday\"],\n  [\"498\u00d79=4482\", \"691\u00d72=1382\"],\n  [\"751\u00d74=3004\", \"764\u00d73=2292\"],\n  [\"209\u00d74=836\", \"835\u00d73=2505\"],\n  [\"495\u00d73=1485\", \"901\u00d76=5406\"],\n  [\"356\u00d76=2136\", \"300\u00d73=900\"],\n  [\"380\u00d72=760\", \"778\u00d76=4668\"],\n  [\"865\u00d77=6055\", \"897\u00d73=2691\"],\n  [\"736\u00d73=2208\", \"322\u00d79=2898\"],\n  [\"719\u00d79=6471\", \"383\u00d72=766\"],\n  [\"948\u00d74=3792\", \"515\u00d72=1030\"],\n  [\"647\u00d74=2588\", \"121\u00d73=363\"],\n  [\"320\u00d73=960\", \"975\u00d75=4875\"],\n  [\"591\u00d78=4728\", \"151\u00d76=906\"],\n  [\"863\u00d74=3452\", \"251\u00d73=753\"],\n  [\"535\u00d78=4280\", \"251\u00d72=502\"],\n  [\"561\u00d77=3927\", \"831\u00d72=1662\"],\n  [\"489\u00d78=3912\", \"370\u00d79=3330\"],\n  [\"469\u00d72=938\", \"802\u00d74=3208\"],\n  [\"457\u00d74=1828\", \"421\u00d73=1263\"],\n  [\"670\u00d72=1340\", \"822\u00d74=3288\"],\n  [\"135\u00d78=1080\", \"177\u00d72=354\"],\n  [\"106\u00d76=636\", \"868\u00d77=6076\"],\n  [\"780\u00d72=1560\", \"183\u00d77=1281\"],\n  [\"116\u00d75=580\", \"763\u00d76=4578\"],\n  [\"863\u00d72=1726\", \"910\u00d79=8190\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every three-digit x one-digit answer cell\n# to the new values, per the commit's regenerated output.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-07-25 Thursday\", \"2024-07-26 Friday\"),\n    @(\"498\u00d79=4482\", \"691\u00d72=1382\"),\n    @(\"751\u00d74=3004\", \"764\u00d73=2292\"),\n    @(\"209\u00d74=836\", \"835\u00d73=2505\"),\n    @(\"495\u00d73=1485\", \"901\u00d76=5406\"),\n    @(\"356\u00d76=2136\", \"300\u00d73=900\"),\n    @(\"380\u00d72=760\", \"778\u00d76=4668\"),\n    @(\"865\u00d77=6055\", \"897\u00d73=2691\"),\n    @(\"736\u00d73=2208\", \"322\u00d79=2898\"),\n    @(\"719\u00d79=6471\", \"383\u00d72=766\"),\n    @(\"948\u00d74=3792\", \"515\u00d72=1030\"),\n    @(\"647\u00d74=2588\", \"121\u00d73=363\"),\n    @(\"320\u00d73=960\", \"975\u00d75=4875\"),\n    @(\"591\u00d78=4728\", \"151\u00d76=906\"),\n    @(\"863\u00d74=3452\", \"251\u00d73=753\"),\n    @(\"535\u00d78=4280\", \"251\u00d72=502\"),\n    @(\"561\u00d77=3927\", \"831\u00d72=1662\"),\n    @(\"489\u00d78=3912\", \"370\u00d79=3330\"),\n    @(\"469\u00d72=938\", \"802\u00d74=3208\"),\n    @(\"457\u00d74=1828\", \"421\u00d73=1263\"),\n    @(\"670\u00d72=1340\", \"822\u00d74=3288\"),\n    @(\"135\u00d78=1080\", \"177\u00d72=354\"),\n    @(\"106\u00d76=636\", \"868\u00d77=6076\"),\n    @(\"780\u00d72=1560\", \"183\u00d77=1281\"),\n    @(\"116\u00d75=580\", \"763\u00d76=4578\"),\n    @(\"863\u00d72=1726\", \"910\u00d79=8190\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $oldText\"\n    }\n}\n"}
